# Corrección año de reporte ventas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder product list text in A2
$ws.Range("A2").Value = "2-Queque,1-Torta,"

# Orders still awaiting pickup -> "Pendiente"
$ws.Range("G2:G6").Value = "Pendiente"

# Orders completed -> "Finalizado"
$ws.Range("G7:G9").Value = "Finalizado"

# Fix the estimated pickup date year/month typo for row 9
$ws.Range("C9").Value = "20-10-2019"
